$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'312.52"
$ws.Range("E2").Value = "'2.03%"

$ws.Range("D3").Value = "'37.49"
$ws.Range("E3").Value = "'-0.02%"

$ws.Range("D4").Value = "'5.136"
$ws.Range("E4").Value = "'0.87%"

$ws.Range("D5").Value = "'0.07857"
$ws.Range("E5").Value = "'1.71%"

$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.910"
$ws.Range("E6").Value = "'0.84%"

$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'8.277"
$ws.Range("E7").Value = "'0.85%"

$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.908"
$ws.Range("E8").Value = "'-9.00%"

$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9175"
$ws.Range("E9").Value = "'0.08%"

$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1186"
$ws.Range("E10").Value = "'-0.23%"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1917"
$ws.Range("E11").Value = "'1.74%"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09042"
$ws.Range("E12").Value = "'3.91%"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03341"
$ws.Range("E13").Value = "'-1.91%"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09586"
$ws.Range("E14").Value = "'-1.22%"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001387"
$ws.Range("E15").Value = "'1.34%"

$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005713"
$ws.Range("E16").Value = "'-3.93%"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.510"
$ws.Range("E17").Value = "'-2.06%"

$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.414"
$ws.Range("E18").Value = "'1.34%"

$ws.Range("D20").Value = "'5.241"
$ws.Range("E20").Value = "'4.43%"

$ws.Range("D21").Value = "'0.1274"
$ws.Range("E21").Value = "'-0.12%"

$ws.Range("D22").Value = "'0.2590"
$ws.Range("E22").Value = "'-0.21%"

$ws.Range("D23").Value = "'0.04368"
$ws.Range("E23").Value = "'0.92%"

$ws.Range("D24").Value = "'0.001249"
$ws.Range("E24").Value = "'2.86%"

$ws.Range("D25").Value = "'0.004688"
$ws.Range("E25").Value = "'3.17%"

$ws.Range("D26").Value = "'0.0001359"
$ws.Range("E26").Value = "'0.52%"

$ws.Range("D27").Value = "'0.0003989"
$ws.Range("E27").Value = "'-98.11%"

$ws.Range("D39").Value = "'0.02309"
$ws.Range("E39").Value = "'4.34%"

$ws.Range("D40").Value = "'0.05071"
$ws.Range("E40").Value = "'3.17%"

$ws.Range("D41").Value = "'0.007468"
$ws.Range("E41").Value = "'-1.12%"

$ws.Range("D42").Value = "'0.009044"
$ws.Range("E42").Value = "'-8.54%"

$ws.Range("D43").Value = "'0.1350"
$ws.Range("E43").Value = "'1.06%"

$ws.Range("D44").Value = "'0.001950"
$ws.Range("E44").Value = "'-5.51%"

$ws.Range("D45").Value = "'0.009394"
$ws.Range("E45").Value = "'6.68%"

$ws.Range("D46").Value = "'0.00006627"
$ws.Range("E46").Value = "'1.06%"

$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.19%"

$ws.Range("D48").Value = "'0.003349"
$ws.Range("E48").Value = "'11.67%"

$ws.Range("D49").Value = "'0.0009997"
$ws.Range("E49").Value = "'-23.23%"

$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.19%"

$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.19%"
